$d = $word.ActiveDocument

# ------------------------------------------------------------------
# Context: paragraphs (1-based) around the "Run 11..15" bridge-run
# notes. Run 11 and Run 12's visible text do not change at all, so
# those paragraphs are left untouched. Run 13 gets "(done)" appended.
# Run 14's old text ("add AR1 on p move north to south") becomes the
# new Run 15 content, and what used to be Run 15's text
# ("switch to logistic-normal for age comp") becomes the new Run 14.
# A brand new Run 16 paragraph is appended with the estimate-sd-scalar
# note, and the _GoBack bookmark (currently sitting inside the old
# Run 12 paragraph) is moved to the very end of the new Run 16 text.
# ------------------------------------------------------------------

# 1) Strip the _GoBack bookmark from its old spot (inside the Run 12
#    paragraph) - it will be re-added at the end of the new Run 16
#    paragraph further down. The bookmark is collapsed (empty), so
#    deleting it does not touch any visible text.
if ($d.Bookmarks.Exists("_GoBack")) {
    $d.Bookmarks.Item("_GoBack").Delete()
}

# 2) Run 13: append "(done)" right after "...at 0.1"
$p13 = $d.Paragraphs.Item(32)
$r13 = $d.Range($p13.Range.Start, $p13.Range.End - 1)
$r13.InsertAfter("(done)")

# 3) Run 14: old text was "Run 14: add AR1 on p move north to south";
#    it becomes "Run 14: switch to logistic-normal for age comp "
$p14 = $d.Paragraphs.Item(34)
$r14 = $d.Range($p14.Range.Start, $p14.Range.End - 1)
$r14.Text = "Run 14: switch to logistic-normal for age comp "

# 4) Run 15: old text was "Run 15: switch to logistic-normal for age
#    comp"; it becomes the "estimate sd scalar..." note.
$p15 = $d.Paragraphs.Item(36)
$r15 = $d.Range($p15.Range.Start, $p15.Range.End - 1)
$r15.Text = "Run 15: estimate sd scalar for aggregate rec cpa indices (doesn" + [char]8217 + "t converge)"

# 5) Add a new Run 16 paragraph after the blank paragraph that follows
#    Run 15.
$pBlank = $d.Paragraphs.Item(37)
$pBlank.Range.InsertParagraphAfter()
$p16 = $d.Paragraphs.Item(38)
$r16 = $d.Range($p16.Range.Start, $p16.Range.End - 1)
$r16.Text = "Run 16: add AR1 on p move north to south"

# 6) Re-add the _GoBack bookmark at the very end of the new Run 16
#    text. A collapsed bookmark placed exactly at a paragraph's last
#    valid offset lands in the wrong spot, so a scratch character is
#    appended first, the bookmark is anchored just before it, and the
#    scratch character is then deleted - leaving the bookmark as the
#    last thing in the paragraph.
$p16b = $d.Paragraphs.Item(38)
$scratchPos = $p16b.Range.End - 1
$d.Range($scratchPos, $scratchPos).InsertAfter("X")

$p16c = $d.Paragraphs.Item(38)
$bmPos = $p16c.Range.End - 2
$d.Bookmarks.Add("_GoBack", $d.Range($bmPos, $bmPos)) | Out-Null

$p16d = $d.Paragraphs.Item(38)
$d.Range($p16d.Range.End - 2, $p16d.Range.End - 1).Delete()

Write-Output "done"
